$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K34").Value = 'Iraq (2), Pakistan (2)'
$ws.Range("K37").Value = 'Nigeria (4), Russia (3), NA (3), Iraq (2), Pakistan (2), Uzbekistan (1)'
$ws.Range("K46").Value = 'Russia (11)'
$ws.Range("K47").Value = 'China (1)'
$ws.Range("K48").Value = 'China (3)'
$ws.Range("K49").Value = 'China (1)'
$ws.Range("K51").Value = 'Russia (1)'
$ws.Range("K52").Value = 'Turkey (1)'
$ws.Range("K58").Value = 'NA (1), Nigeria (1)'
$ws.Range("K70").Value = 'Russia (4)'
$ws.Range("K71").Value = 'Russia (1)'
$ws.Range("K75").Value = 'China (2)'
$ws.Range("K76").Value = 'Nigeria (2)'
$ws.Range("K77").Value = 'Pakistan (1)'
$ws.Range("K99").Value = 'Russia (10)'
$ws.Range("K101").Value = 'Kazakhstan (5)'
$ws.Range("K102").Value = 'Kazakhstan (11)'
$ws.Range("K103").Value = 'Ghana (3)'
$ws.Range("K106").Value = 'Turkey (9)'
$ws.Range("K110").Value = 'Sudan (9)'
$ws.Range("K111").Value = 'Sudan (2)'
$ws.Range("K121").Value = 'Ghana (3)'
$ws.Range("K123").Value = 'Sudan (3)'
$ws.Range("K126").Value = 'Kazakhstan (3)'
$ws.Range("K127").Value = 'Turkey (4)'
$ws.Range("K132").Value = 'Kosovo (40)'
$ws.Range("K135").Value = 'United Arab Emirates (15)'
$ws.Range("K154").Value = 'Russia (6)'
$ws.Range("K155").Value = 'Kazakhstan (18)'
$ws.Range("K156").Value = 'South Africa (9), Greece (6), Russia (6), Senegal (6), Republic of the Congo (3), Turkey (3), Oman (3), China (3), Mauritania (3)'
$ws.Range("K162").Value = 'Yugoslavia (1)'
$ws.Range("K202").Value = 'Algeria (2)'
$ws.Range("K205").Value = 'Russia (12)'
$ws.Range("K207").Value = 'Russia (42), Tajikistan (8), Kazakhstan (5), Turkmenistan (1), Uzbekistan (1), Bulgaria (1)'
$ws.Range("K211").Value = 'Russia (23), Tajikistan (8), NA (8), Kazakhstan (5), Turkmenistan (1), Uzbekistan (1), Bulgaria (1)'
$ws.Range("K212").Value = 'Russia (24), Turkmenistan (6), Uganda (3), Kazakhstan (3), Nigeria (3), Uzbekistan (3), Tajikistan (3)'
$ws.Range("K224").Value = 'Russia (2)'
$ws.Range("K225").Value = 'Russia (9), Tajikistan (2), Bulgaria (1), Turkmenistan (1), Uzbekistan (1)'
$ws.Range("K226").Value = 'Tajikistan (1), Russia (1)'
$ws.Range("K234").Value = 'China (7)'
$ws.Range("K239").Value = 'Russia (3)'
$ws.Range("K241").Value = 'Uganda (2)'
$ws.Range("K243").Value = 'Kenya (8)'
$ws.Range("K253").Value = 'Turkey (2)'
$ws.Range("K260").Value = 'Greece (2)'
$ws.Range("K264").Value = 'Greece (1)'
$ws.Range("K269").Value = 'Russia (6)'
$ws.Range("K271").Value = 'Tajikistan (3)'
$ws.Range("K272").Value = 'Tajikistan (9)'
$ws.Range("K281").Value = 'India (165)'
$ws.Range("K283").Value = 'China (17)'
$ws.Range("K286").Value = 'Nigeria (2), Pakistan (1)'
$ws.Range("K289").Value = 'Russia (2), Tajikistan (1)'
$ws.Range("K294").Value = 'Turkey (2), Syria (2)'
$ws.Range("K297").Value = 'USA (5)'
$ws.Range("K298").Value = 'China (2), NA (2)'
$ws.Range("K299").Value = 'China (2)'
$ws.Range("K306").Value = 'Russia (3)'
$ws.Range("K307").Value = 'Uganda (3)'
$ws.Range("K308").Value = 'Turkey (3)'
$ws.Range("K309").Value = 'India (24)'
$ws.Range("K311").Value = 'India (21)'
$ws.Range("K312").Value = 'Russia (14)'
$ws.Range("K314").Value = 'Saudi Arabia (3)'
$ws.Range("K315").Value = 'Turkey (11)'
$ws.Range("K321").Value = 'Mali (3)'
